$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Carrera" column (B) changes from "Negocios" to "Sistemas" for both data rows
$ws.Range("B2").Value = "Sistemas"
$ws.Range("B3").Value = "Sistemas"
